$d = $word.ActiveDocument

# Paragraph 3 ("3- Poderia ser aplicado...") currently carries the _GoBack
# bookmark wrapping its whole run. We need to:
#   1) insert a brand-new paragraph 4 ("4- Seria uma equipe...") right after it
#   2) move the _GoBack bookmark so it collapses at the end of the new
#      paragraph 4's run (after the text, before the paragraph mark)

$p3 = $d.Paragraphs.Item(3)

# Create the new (empty) paragraph right after paragraph 3.
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs.Item(4)

# Re-anchor the _GoBack bookmark onto the brand-new, still-empty paragraph
# *before* typing its text. Bookmark names are unique, so adding a bookmark
# named "_GoBack" here automatically removes the old one that was wrapping
# paragraph 3. Because the paragraph is empty when the bookmark is added,
# the collapsed bookmark naturally ends up wrapping around (and then,
# once text is inserted, trailing) the run that gets typed into it.
$d.Bookmarks.Add("_GoBack", $p4.Range.Duplicate)

# Now fill in paragraph 4's text; it lands between bookmarkStart/bookmarkEnd,
# which is exactly the target shape (run, then bookmarkEnd).
$p4.Range.Text = "4- Seria uma equipe de seis a oito pessoas. Como é um ciclo separaríamos em partes o todo o projeto, após o levantamento de requisitos e ter organizado o projeto para ser entregue em pedaços iríamos fazer cada pedaço juntos cada um cuidando de uma parte e caso precisar de ajuda outro membro do grupo auxilia."
